# Scheduled market-data refresh: overwrite currentAveragePrice* / Leve profit
# columns (H:N) for the leves whose underlying item prices moved, across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. Plain numeric overwrites (no
# formulas are stored in these cells) - a couple of cells whose profit
# figure is no longer computable are cleared outright rather than zeroed.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1389.875
$ws.Range("I96").Value = 1918.5714
$ws.Range("J96").Value = 978.6667
$ws.Range("K96").Value = 5755.7142
$ws.Range("L96").Value = 2936.0001
$ws.Range("M96").Value = -4382.7142
$ws.Range("N96").Value = -5682.0001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 112396.78
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 112396.78
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 337190.34
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -338182.34

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1740.4166
$ws.Range("I106").Value = 1737.2222
$ws.Range("J106").Value = 1750
$ws.Range("K106").Value = 1737.2222
$ws.Range("L106").Value = 1750
$ws.Range("M106").Value = -1106.2222
$ws.Range("N106").Value = -3012

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1985569
$ws.Range("I132").Value = 2416404.2
$ws.Range("J132").Value = 3727.2
$ws.Range("K132").Value = 7249212.600000001
$ws.Range("L132").Value = 11181.6
$ws.Range("M132").Value = -7246682.600000001
$ws.Range("N132").Value = -16241.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1353.8
$ws.Range("I2").Value = 882.7692
$ws.Range("J2").Value = 2228.5715
$ws.Range("K2").Value = 882.7692
$ws.Range("L2").Value = 2228.5715
$ws.Range("M2").Value = -769.7692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5117.018
$ws.Range("I32").Value = 3335.46
$ws.Range("J32").Value = 19963.334
$ws.Range("K32").Value = 3335.46
$ws.Range("L32").Value = 19963.334
$ws.Range("M32").Value = -3048.46
$ws.Range("N32").Value = -20537.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1353.8
$ws.Range("I116").Value = 882.7692
$ws.Range("J116").Value = 2228.5715
$ws.Range("K116").Value = 882.7692
$ws.Range("L116").Value = 2228.5715
$ws.Range("M116").Value = 1411.2308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3403.8572
$ws.Range("I132").Value = 1750
$ws.Range("J132").Value = 4065.4
$ws.Range("K132").Value = 5250
$ws.Range("L132").Value = 12196.2
$ws.Range("M132").Value = -2720
$ws.Range("N132").Value = -17256.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1353.8
$ws.Range("I3").Value = 882.7692
$ws.Range("J3").Value = 2228.5715
$ws.Range("K3").Value = 882.7692
$ws.Range("L3").Value = 2228.5715
$ws.Range("M3").Value = -768.7692

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1163.3077
$ws.Range("I107").Value = 1145.625
$ws.Range("J107").Value = 1191.6
$ws.Range("K107").Value = 1145.625
$ws.Range("L107").Value = 1191.6
$ws.Range("M107").Value = 774.375
$ws.Range("N107").Value = -5031.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2633.111
$ws.Range("I58").Value = 1900
$ws.Range("J58").Value = 2724.75
$ws.Range("K58").Value = 1900
$ws.Range("L58").Value = 2724.75
$ws.Range("M58").Value = -1697
$ws.Range("N58").Value = -3130.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2390.8333
$ws.Range("I134").Value = 1916.6666
$ws.Range("J134").Value = 3181.111
$ws.Range("K134").Value = 5749.9998
$ws.Range("L134").Value = 9543.332999999999
$ws.Range("M134").Value = -3214.9998
$ws.Range("N134").Value = -14613.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2633.111
$ws.Range("I136").Value = 1900
$ws.Range("J136").Value = 2724.75
$ws.Range("K136").Value = 5700
$ws.Range("L136").Value = 8174.25
$ws.Range("M136").Value = -3150
$ws.Range("N136").Value = -13274.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 6000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 6000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18228

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2372.6
$ws.Range("I80").Value = 2150
$ws.Range("J80").Value = 2428.25
$ws.Range("K80").Value = 2150
$ws.Range("L80").Value = 2428.25
$ws.Range("M80").Value = -1152
$ws.Range("N80").Value = -4424.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2372.6
$ws.Range("I83").Value = 2150
$ws.Range("J83").Value = 2428.25
$ws.Range("K83").Value = 10750
$ws.Range("L83").Value = 12141.25
$ws.Range("M83").Value = -5758
$ws.Range("N83").Value = -22125.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1176.2222
$ws.Range("I107").Value = 985.75
$ws.Range("J107").Value = 2700
$ws.Range("K107").Value = 985.75
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = 934.25
$ws.Range("N107").Value = -6540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5840.5
$ws.Range("I113").Value = 7688.75
$ws.Range("J113").Value = 2144
$ws.Range("K113").Value = 7688.75
$ws.Range("L113").Value = 2144
$ws.Range("M113").Value = -5518.75
$ws.Range("N113").Value = -6484

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8186.75
$ws.Range("I132").Value = 11929.363
$ws.Range("J132").Value = 3612.4443
$ws.Range("K132").Value = 35788.089
$ws.Range("L132").Value = 10837.3329
$ws.Range("M132").Value = -33258.089
$ws.Range("N132").Value = -15897.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1920.7931
$ws.Range("I7").Value = 1915.3077
$ws.Range("J7").Value = 1968.3334
$ws.Range("K7").Value = 1915.3077
$ws.Range("L7").Value = 1968.3334
$ws.Range("M7").Value = -1803.3077

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 392.41464
$ws.Range("I22").Value = 384.97223
$ws.Range("J22").Value = 446
$ws.Range("K22").Value = 384.97223
$ws.Range("L22").Value = 446
$ws.Range("M22").Value = -89.97223000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 392.41464
$ws.Range("I27").Value = 384.97223
$ws.Range("J27").Value = 446
$ws.Range("K27").Value = 384.97223
$ws.Range("L27").Value = 446
$ws.Range("M27").Value = -277.97223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1920.7931
$ws.Range("I126").Value = 1915.3077
$ws.Range("J126").Value = 1968.3334
$ws.Range("K126").Value = 5745.9231
$ws.Range("L126").Value = 5905.0002
$ws.Range("M126").Value = -3275.9231

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8806.267
$ws.Range("I81").Value = 12844.444
$ws.Range("J81").Value = 2749
$ws.Range("K81").Value = 25688.888
$ws.Range("L81").Value = 5498
$ws.Range("M81").Value = -24627.888
$ws.Range("N81").Value = -7620

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 8806.267
$ws.Range("I84").Value = 12844.444
$ws.Range("J84").Value = 2749
$ws.Range("K84").Value = 128444.44
$ws.Range("L84").Value = 27490
$ws.Range("M84").Value = -123140.44
$ws.Range("N84").Value = -38098

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5347995
$ws.Range("I100").Value = 5347995
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 10695990
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -10695449

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3700.182
$ws.Range("I107").Value = 679.7273
$ws.Range("J107").Value = 9741.091
$ws.Range("K107").Value = 2039.1819
$ws.Range("L107").Value = 29223.273
$ws.Range("M107").Value = -119.1819
$ws.Range("N107").Value = -33063.273
